# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the first data row (65c295ce-4e6e-40bb-b932-1f3e43f25f31) on the
# zh-cn and de-de language sheets, reflecting the latest handback report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-10-20 08:56:10"
$wsZhCn.Range("K2").Value = "2016-10-20 08:56:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-10-20 08:56:22"
$wsDeDe.Range("K2").Value = "2016-10-20 08:57:13"

# The Overview sheet's "Latest HO Xliff Generate Date" column mirrors the
# de-de sheet's "Correspond Handoff Datetime" for each file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-20 08:56:22"
